$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before N ("Total salary") for the new "Bonus reward" column.
$ws.Columns("N:N").Insert() | Out-Null

# Set the header text for the newly inserted column.
$ws.Range("N6").Value = "Bonus reward"

# Re-apply the autofilter so its range grows to include the new column.
$ws.AutoFilterMode = $false
$ws.Range("A6:Q6").AutoFilter() | Out-Null

# The hidden _FilterDatabase defined name also needs to track the new range.
$fdb = $wb.Names.Item(1)
$fdb.RefersTo = "=Salary!`$A`$6:`$Q`$6"

# Update the view: clear the old scroll position and move the selection.
$ws.Range("F15").Select() | Out-Null
